$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 18:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 449555
$ws.Range("C4").Value = 14628
$ws.Range("D4").Value = 24562
$ws.Range("E4").Value = 409167
$ws.Range("F4").Value = 9704
$ws.Range("G4").Value = 1038
$ws.Range("H4").Value = 15826

# --- Row 12: Turquia ---
$ws.Range("B12").Value = 42282
$ws.Range("C12").Value = 4056
$ws.Range("D12").Value = 2142
$ws.Range("E12").Value = 39232
$ws.Range("F12").Value = 1552
$ws.Range("G12").Value = 96
$ws.Range("H12").Value = 908

# --- Row 14: Suiza ---
$ws.Range("B14").Value = 24046
$ws.Range("C14").Value = 766
$ws.Range("D14").Value = 10600
$ws.Range("E14").Value = 12498
$ws.Range("G14").Value = 53
$ws.Range("H14").Value = 948

# --- Row 24: India ---
$ws.Range("B24").Value = 6653
$ws.Range("C24").Value = 737
$ws.Range("E24").Value = 5856

# --- Row 50: Colombia ---
$ws.Range("F50").Value = 85

# --- Reorder "Monaco" / "Guayana Francesa" and update their stats ---
# Row 133 was Guayana Francesa, Row 134 was Monaco.
# After the edit, row 133 becomes Monaco (with new stats) and row 134 becomes
# Guayana Francesa (keeping its former, unchanged stats).
$ws.Range("A133").Value = "Monaco"
$ws.Range("B133").Value = 84
$ws.Range("C133").Value = 3
$ws.Range("D133").Value = 5
$ws.Range("E133").Value = 78
$ws.Range("F133").Value = 4
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 1

$ws.Range("A134").Value = "Guayana Francesa"
$ws.Range("B134").Value = 83
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 43
$ws.Range("E134").Value = 40
$ws.Range("F134").Value = 1
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0
